# Actualización automática del tracker
# - Rellena resultado/profit de dos partidos ya resueltos (filas 101 y 107)
# - Añade los nuevos partidos pendientes de resultado (filas 124-130)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, [string]$text) {
    # Forzar que el valor se guarde como texto, evitando que Excel lo
    # reinterprete como fecha/numero (p.ej. "2025-08-08"), y despues
    # limpiar el formato para que la celda quede sin estilo adicional.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

function Set-BlankCell($cell, $templateCell) {
    # Deja la celda presente pero vacia, igual que las celdas
    # "resultado"/"profit" de los partidos que aun no se han jugado.
    $templateCell.Copy($cell)
}

# ---------------------------------------------------------------
# Resultados que ya se conocen: completar resultado (G) y profit (H)
# ---------------------------------------------------------------
Set-TextValue $ws.Range("G101") "Fallo"
$ws.Range("H101").Value = -1

Set-TextValue $ws.Range("G107") "Fallo"
$ws.Range("H107").Value = -1

# ---------------------------------------------------------------
# Nuevos partidos añadidos al tracker (sin resultado todavía)
# ---------------------------------------------------------------
$newMatches = @(
    @{ Row=124; EventId=14358000; Fecha="2025-08-08"; JugadorA="Francisco Comesaña"; JugadorB="Jaume Munar"; Pronostico="Gana Francisco Comesaña"; Cuota=3 },
    @{ Row=125; EventId=14357978; Fecha="2025-08-08"; JugadorA="Hamad Medjedovic"; JugadorB="Aleksandar Kovacevic"; Pronostico="Gana Aleksandar Kovacevic"; Cuota=1.91 },
    @{ Row=126; EventId=14357962; Fecha="2025-08-08"; JugadorA="Daniel Altmaier"; JugadorB="Roberto Bautista Agut"; Pronostico="Gana Daniel Altmaier"; Cuota=2.5 },
    @{ Row=127; EventId=14379989; Fecha="2025-08-09"; JugadorA="Aleksandar Vukic"; JugadorB="Nishesh Basavareddy"; Pronostico="Gana Nishesh Basavareddy"; Cuota=2.2 },
    @{ Row=128; EventId=14366982; Fecha="2025-08-08"; JugadorA="Lin Zhu"; JugadorB="Lucia Bronzetti"; Pronostico="Gana Lucia Bronzetti"; Cuota=2.3 },
    @{ Row=129; EventId=14311733; Fecha="2025-08-08"; JugadorA="Federico Cinà"; JugadorB="Ugo Blanchet"; Pronostico="Gana Ugo Blanchet"; Cuota=2.38 },
    @{ Row=130; EventId=14311065; Fecha="2025-08-08"; JugadorA="Carlos Taberner"; JugadorB="Gonzalo Bueno"; Pronostico="Gana Gonzalo Bueno"; Cuota=3 }
)

foreach ($m in $newMatches) {
    $r = $m.Row

    $ws.Range("A$r").Value = $m.EventId
    Set-TextValue $ws.Range("B$r") $m.Fecha
    $ws.Range("C$r").Value = $m.JugadorA
    $ws.Range("D$r").Value = $m.JugadorB
    $ws.Range("E$r").Value = $m.Pronostico
    $ws.Range("F$r").Value = $m.Cuota

    Set-BlankCell $ws.Range("G$r") $ws.Range("G123")
    Set-BlankCell $ws.Range("H$r") $ws.Range("H123")
}
